# =====================================================================
# Rebuild "B0CTTY491F_sales_po_comparison.xlsx" forecast workbook:
#  - Rename Sheet1 -> "Sales vs PO"; insert "Order Week" column (old ds),
#    shift ds forward one order-cycle, zero out PO_Requested_Qty.
#  - Add "Weekly Growth" sheet: ds/PO_Requested_Qty/Growth% for weeks
#    that actually had PO activity.
#  - Add "Volume Insights" sheet: aggregate PO stats.
#  - Add "Prediction Info" sheet: next-week PO forecast.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Sales vs PO"
# ---------------------------------------------------------------------
$wsSales = $wb.ActiveSheet
$wsSales.Name = "Sales vs PO"

# New data table: col1=ds(shifted +6d), col2=y(unchanged),
#                 col3=Order Week(original ds), col4=PO_Requested_Qty(0)
$salesData = @(
    @(45340,1,45334,0),
    @(45347,0,45341,0),
    @(45354,3,45348,0),
    @(45361,41,45355,0),
    @(45368,37,45362,0),
    @(45375,19,45369,0),
    @(45382,19,45376,0),
    @(45389,3,45383,0),
    @(45396,3,45390,0),
    @(45403,3,45397,0),
    @(45410,1,45404,0),
    @(45417,0,45411,0),
    @(45424,0,45418,0),
    @(45431,5,45425,0),
    @(45438,4,45432,0),
    @(45445,12,45439,0),
    @(45452,5,45446,0),
    @(45459,4,45453,0),
    @(45466,0,45460,0),
    @(45473,1,45467,0),
    @(45480,3,45474,0),
    @(45487,7,45481,0),
    @(45494,4,45488,0),
    @(45501,2,45495,0),
    @(45508,14,45502,0),
    @(45515,18,45509,0),
    @(45522,16,45516,0),
    @(45529,70,45523,0),
    @(45536,17,45530,0),
    @(45543,7,45537,0),
    @(45550,4,45544,0),
    @(45557,9,45551,0),
    @(45564,23,45558,0),
    @(45571,13,45565,0),
    @(45578,23,45572,0),
    @(45585,11,45579,0),
    @(45592,7,45586,0),
    @(45599,24,45593,0),
    @(45606,29,45600,0),
    @(45613,34,45607,0),
    @(45620,33,45614,0),
    @(45627,15,45621,0),
    @(45634,21,45628,0),
    @(45641,17,45635,0),
    @(45648,3,45642,0),
    @(45655,3,45649,0)
)

# Give column C the same date number-format as column A (rows 2-47),
# and give the new header cells the same style as the existing headers.
$wsSales.Range("A2").Copy()
$wsSales.Range("C2:C47").PasteSpecial(-4122)
$wsSales.Range("A1").Copy()
$wsSales.Range("C1:D1").PasteSpecial(-4122)

$wsSales.Range("C1").Value = "Order Week"
$wsSales.Range("D1").Value = "PO_Requested_Qty"

$r = 2
foreach ($row in $salesData) {
    $wsSales.Cells.Item($r, 1).Value = $row[0]
    $wsSales.Cells.Item($r, 2).Value = $row[1]
    $wsSales.Cells.Item($r, 3).Value = $row[2]
    $wsSales.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# ---------------------------------------------------------------------
$wsGrowth = $wb.Worksheets.Add($null, $wsSales)
$wsGrowth.Name = "Weekly Growth"
$wsGrowth.PageSetup.LeftMargin = 54
$wsGrowth.PageSetup.RightMargin = 54
$wsGrowth.PageSetup.TopMargin = 72
$wsGrowth.PageSetup.BottomMargin = 72
$wsGrowth.PageSetup.HeaderMargin = 36
$wsGrowth.PageSetup.FooterMargin = 36

$wsSales.Range("A1").Copy()
$wsGrowth.Range("A1:C1").PasteSpecial(-4122)
$wsGrowth.Range("A1").Value = "ds"
$wsGrowth.Range("B1").Value = "PO_Requested_Qty"
$wsGrowth.Range("C1").Value = "Growth%"

$growthData = @(
    @(45341,336,0),
    @(45355,32,-90.47619047619048),
    @(45418,192,500.0),
    @(45425,16,-91.66666666666666),
    @(45439,96,500.0),
    @(45453,16,-83.33333333333334),
    @(45516,16,0.0),
    @(45523,32,100.0),
    @(45530,48,50.0),
    @(45537,64,33.33333333333333),
    @(45558,96,50.0),
    @(45572,16,-83.33333333333334),
    @(45579,128,700.0)
)

$wsSales.Range("A2").Copy()
$wsGrowth.Range("A2:A" + (1 + $growthData.Count)).PasteSpecial(-4122)

$r = 2
foreach ($row in $growthData) {
    $wsGrowth.Cells.Item($r, 1).Value = $row[0]
    $wsGrowth.Cells.Item($r, 2).Value = $row[1]
    $wsGrowth.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ---------------------------------------------------------------------
$wsVolume = $wb.Worksheets.Add($null, $wsGrowth)
$wsVolume.Name = "Volume Insights"
$wsVolume.PageSetup.LeftMargin = 54
$wsVolume.PageSetup.RightMargin = 54
$wsVolume.PageSetup.TopMargin = 72
$wsVolume.PageSetup.BottomMargin = 72
$wsVolume.PageSetup.HeaderMargin = 36
$wsVolume.PageSetup.FooterMargin = 36

$wsSales.Range("A1").Copy()
$wsVolume.Range("A1:D1").PasteSpecial(-4122)
$wsVolume.Range("A1").Value = "Total_PO_Quantity"
$wsVolume.Range("B1").Value = "Average_PO_Quantity"
$wsVolume.Range("C1").Value = "Max_PO_Quantity"
$wsVolume.Range("D1").Value = "Min_PO_Quantity"

$wsVolume.Range("A2").Value = 1088
$wsVolume.Range("B2").Value = 83.69230769230769
$wsVolume.Range("C2").Value = 336
$wsVolume.Range("D2").Value = 16

# ---------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ---------------------------------------------------------------------
$wsPredict = $wb.Worksheets.Add($null, $wsVolume)
$wsPredict.Name = "Prediction Info"
$wsPredict.PageSetup.LeftMargin = 54
$wsPredict.PageSetup.RightMargin = 54
$wsPredict.PageSetup.TopMargin = 72
$wsPredict.PageSetup.BottomMargin = 72
$wsPredict.PageSetup.HeaderMargin = 36
$wsPredict.PageSetup.FooterMargin = 36

$wsSales.Range("A1").Copy()
$wsPredict.Range("A1").PasteSpecial(-4122)
$wsPredict.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"

$wsPredict.Range("A2").Value = 20.30769230769228

# Keep "Sales vs PO" as the active/selected tab (matches source activeTab="0").
$wsSales.Activate()

Write-Output "Workbook rebuilt: 4 sheets (Sales vs PO, Weekly Growth, Volume Insights, Prediction Info)."
